$d = $word.ActiveDocument

# --- 1) Merge the split "{#is_" + "title}" runs into a single run "{#is_title}" ---
# (Find/Replace across the run boundary makes Word coalesce the runs into one,
#  which also drops the stray <w:proofErr w:type="gramStart"/> marker that had
#  been forcing the split.)
$d.Content.Find.Execute("{#is_title}", $false, $false, $false, $false, $false, $true, 1, $false, "{#is_title}", 2) | Out-Null

# --- 2) Merge the split "{" + "message}" runs into a single run "{message}" ---
# (same coalescing trick; removes the matching <w:proofErr w:type="gramEnd"/>)
$d.Content.Find.Execute("{message}", $false, $false, $false, $false, $false, $true, 1, $false, "{message}", 2) | Out-Null

# --- 3) Remove the two whole paragraphs that held the bodycenter / ---
# ---    bodyboldcenter template blocks (and the bookmark inside them) ---
$searchRange = $d.Content
$searchRange.Find.Execute("{#is_bodycenter}") | Out-Null
$pBodyCenter = $searchRange.Paragraphs(1)
$pBodyBoldCenter = $pBodyCenter.Next()

$delRange = $d.Range($pBodyCenter.Range.Start, $pBodyBoldCenter.Range.End)
$delRange.Delete()
